$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the decimal number format for the "Difference in Incident" column (E2:E13)
# and bring B12:D12 along to the same "0.0" number format.
$ws.Range("E2:E13").NumberFormat = "0.0"
$ws.Range("B12:D12").NumberFormat = "0.0"

# Correct the 2020 Safety Incidents value for November (row 12)
$ws.Range("C12").Value = 32

# Re-apply the difference formula as a fill-down so Excel stores it as a shared formula
$ws.Range("E2").Formula = "=B2-C2"
$ws.Range("E2").AutoFill($ws.Range("E2:E13"), 0)

# Update the active selection to match the saved state
$ws.Range("B13").Select()
